$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the now-unused DATA and "Adjusted Expenditure" worksheets; only
# MAIN should remain in the workbook.
$wb.Worksheets("DATA").Delete() | Out-Null
$wb.Worksheets("Adjusted Expenditure").Delete() | Out-Null
